# Generate Report for Handoff
#
# The 091725da-5e97-4558-9d5f-02f3c42e4e54 file finished its localization
# handoff run: its Priority flips from "low" to "ht" (matching the other
# rows) and its Latest Handoff Datetime / Latest HO Xliff Generate Date
# timestamps advance to reflect the new handoff pass, on both the
# per-locale sheets (zh-cn, de-de) and the Overview roll-up sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: Latest HO Xliff Generate Date for 091725da.md (row 4) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-28 20:31:48"

# --- zh-cn sheet: rows 4-7 share Priority "low" -> "ht"; Latest Handoff ---
# --- Datetime for these rows moves from 20:31:28 to 20:31:43          ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("H4").Value = "2016-08-28 20:31:43"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("H5").Value = "2016-08-28 20:31:43"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("H6").Value = "2016-08-28 20:31:43"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("H7").Value = "2016-08-28 20:31:43"

# --- de-de sheet: same Priority flip; Latest Handoff Datetime moves     ---
# --- from 20:31:32 to 20:31:48 (shares the string with Overview!G4)    ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("H4").Value = "2016-08-28 20:31:48"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("H5").Value = "2016-08-28 20:31:48"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("H6").Value = "2016-08-28 20:31:48"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("H7").Value = "2016-08-28 20:31:48"
